$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Round the coordinate values in Q2 and R2
$ws.Range("Q2").Value = 789285
$ws.Range("R2").Value = 7305370

# Remove the Starttid (Z2) and Sluttid (AB2) values entirely
$ws.Range("Z2").ClearContents()
$ws.Range("AB2").ClearContents()
